$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FPCbS")

# Column widths: A widens (no longer auto-fit), new column B gets an explicit width
$ws.Columns.Item(1).ColumnWidth = 25.85546875
$ws.Columns.Item(2).ColumnWidth = 14.42578125

# Row 1 header: A1 becomes an empty bold cell, B1 becomes the full unit label,
# both wrap text; row grows tall enough to show the wrapped text.
$ws.Rows.Item(1).RowHeight = 75

$ws.Cells.Item(1, 1).Value = $null
$ws.Cells.Item(1, 1).Font.Bold = $true
$ws.Cells.Item(1, 1).WrapText = $true

$ws.Cells.Item(1, 2).Value = "FP (flexibility points/MW)"
$ws.Cells.Item(1, 2).HorizontalAlignment = -4152
$ws.Cells.Item(1, 2).WrapText = $true

$ws.Range("B1").Select()

# New rows for additional flexibility-consuming sources
$ws.Cells.Item(15, 1).Value = "crude oil"
$ws.Cells.Item(15, 2).Value = 0

$ws.Cells.Item(16, 1).Value = "heavy or residual fuel oil"
$ws.Cells.Item(16, 2).Value = 0

$ws.Cells.Item(17, 1).Value = "municipal solid waste"
$ws.Cells.Item(17, 2).Value = 0
